# Updates gh-pages output for 上海-漫展信息.xlsx (generated at 456a3b4)
#
# 1) Sheet "展览": insert a new row for the "上海·LOVELIVE ONLY" event
#    (2024-08-03), shifting the existing 2024-08-10.. rows down by one,
#    then re-sequence the column-A index numbers.
# 2) Bump the "想去人数" (column F) counters across all four sheets.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Insert the new row on "展览" at row 42 (before "2024-08-10" row)
# ---------------------------------------------------------------------
$wsExpo = $wb.Worksheets.Item("展览")
$wsExpo.Rows.Item(42).Insert()

# Column B holds a literal text date like "2024-08-03" -- force text
# formatting first so Excel doesn't reinterpret it as a real date serial.
$bCell = $wsExpo.Cells.Item(42, 2)
$bCell.NumberFormat = "@"
$bCell.Value = "2024-08-03"

$wsExpo.Cells.Item(42, 3).Value = "上海·LOVELIVE ONLY"
$wsExpo.Cells.Item(42, 4).Value = "海潮路133号B1 JUMP工坊"
$wsExpo.Cells.Item(42, 5).Value = "2024.08.03 14:00-08.03 19:00"
$wsExpo.Cells.Item(42, 6).Value = 0
$wsExpo.Cells.Item(42, 7).Value = 60
$wsExpo.Cells.Item(42, 8).Value = "https://show.bilibili.com/platform/detail.html?id=86711"
$wsExpo.Cells.Item(42, 9).Value = "//i2.hdslb.com/bfs/openplatform/202405/bllJHQFL1716983812432.jpeg"

# Match the formatting (bold, centered, bordered) already used by the
# other column-A index cells.
$aCell = $wsExpo.Cells.Item(42, 1)
$aCell.Value = 41
$aCell.Font.Bold = $true
$aCell.HorizontalAlignment = -4108
$aCell.VerticalAlignment = -4160
$aCell.Borders.LineStyle = 1

# The row-index column (A) is a plain running count (row - 1); restore
# that sequence for every row pushed down by the insert.
for ($r = 43; $r -le 47; $r++) {
    $wsExpo.Cells.Item($r, 1).Value = $r - 1
}

# ---------------------------------------------------------------------
# 2. Bump the "想去人数" (column F) values
# ---------------------------------------------------------------------
$updates = @{
    "展览" = @(
        @(2,1722), @(3,10055), @(8,1583), @(10,362), @(14,472),
        @(15,1166), @(23,94), @(25,683), @(26,17), @(27,31),
        @(29,218), @(31,331), @(32,213), @(35,547), @(37,522),
        @(38,1267), @(40,368), @(41,324)
    )
    "演出" = @(
        @(4,42), @(10,5), @(19,539), @(21,318), @(22,682), @(23,69)
    )
    "本地生活" = @(
        @(4,800), @(5,186), @(6,2494), @(7,4022), @(10,266), @(11,175)
    )
    "全部类型" = @(
        @(2,1722), @(3,800), @(4,10055), @(5,186), @(7,4022),
        @(8,50), @(9,266), @(10,266), @(12,1583), @(14,362),
        @(17,5), @(18,1166), @(29,318), @(31,683), @(32,69),
        @(35,331), @(37,363), @(39,547), @(42,522), @(44,368), @(46,324)
    )
}

foreach ($sheetName in $updates.Keys) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($pair in $updates[$sheetName]) {
        $row = $pair[0]
        $newVal = $pair[1]
        $ws.Cells.Item($row, 6).Value = $newVal
    }
}
